# Update the "取得日時" (retrieved at) timestamp in column A for rows 2-16
# on the "ランサーズ" sheet to reflect the new append time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-10 18:33:32"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
